$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text (preserves leading zeros, dotted
# version numbers like "27.676.39", and numeric-looking strings like "0.9979")
# without leaving a residual number-format/style change on the cell.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.676.39"
Set-TextValue "E2" "  +0.63%  "

Set-TextValue "D3" "1.777.40"
Set-TextValue "E3" "  +1.66%  "

Set-TextValue "D4" "0.9979"
Set-TextValue "E4" "  -0.54%  "

Set-TextValue "D5" "325.55"
Set-TextValue "E5" "  +0.44%  "

Set-TextValue "D6" "0.9975"
Set-TextValue "E6" "  -0.45%  "

Set-TextValue "D7" "0.4596"
Set-TextValue "E7" "  +4.10%  "

Set-TextValue "D8" "0.3578"
Set-TextValue "E8" "  -0.82%  "

Set-TextValue "D9" "0.07484"
Set-TextValue "E9" "  +0.21%  "

Set-TextValue "D10" "41.76"
Set-TextValue "E10" "  -1.16%  "

Set-TextValue "D11" "1.101"
Set-TextValue "E11" "  +0.35%  "

Set-TextValue "D12" "0.9973"
Set-TextValue "E12" "  -0.56%  "

Set-TextValue "D13" "20.90"
Set-TextValue "E13" "  +1.49%  "

Set-TextValue "D14" "6.040"
Set-TextValue "E14" "  +0.27%  "

Set-TextValue "D15" "7.215"
Set-TextValue "E15" "  +1.16%  "

Set-TextValue "D16" "1.768.64"
Set-TextValue "E16" "  +1.00%  "

Set-TextValue "D17" "93.84"
Set-TextValue "E17" "  +1.86%  "

Set-TextValue "D18" "0.00001058"
Set-TextValue "E18" "  +0.09%  "

Set-TextValue "D19" "0.06445"
Set-TextValue "E19" "  +0.57%  "

Set-TextValue "D20" "0.9972"
Set-TextValue "E20" "  -0.49%  "

Set-TextValue "D21" "17.12"
Set-TextValue "E21" "  +1.72%  "

Set-TextValue "D22" "5.784"
Set-TextValue "E22" "  -1.14%  "

Set-TextValue "D23" "27.740.64"
Set-TextValue "E23" "  +0.67%  "

Set-TextValue "D24" "11.28"
Set-TextValue "E24" "  +1.02%  "

Set-TextValue "D25" "2.074"
Set-TextValue "E25" "  -1.13%  "

Set-TextValue "D26" "165.03"
Set-TextValue "E26" "  +2.29%  "

Set-TextValue "E27" "  -0.54%  "

Set-TextValue "D28" "1.968.30"
Set-TextValue "E28" "  +0.96%  "

Set-TextValue "D29" "2.176"
Set-TextValue "E29" "  +3.49%  "

Set-TextValue "D30" "125.83"
Set-TextValue "E30" "  +1.08%  "

Set-TextValue "D31" "1.095"
Set-TextValue "E31" "  +1.78%  "

Set-TextValue "D32" "0.09219"
Set-TextValue "E32" "  +3.03%  "

Set-TextValue "D33" "3.666"
Set-TextValue "E33" "  +0.41%  "

Set-TextValue "D34" "5.554"
Set-TextValue "E34" "  +0.71%  "

Set-TextValue "D35" "11.85"
Set-TextValue "E35" "  -0.94%  "

Set-TextValue "D36" "0.02292"
Set-TextValue "E36" "  -0.79%  "

Set-TextValue "D37" "0.06120"
Set-TextValue "E37" "  +2.57%  "

Set-TextValue "D38" "0.2095"
Set-TextValue "E38" "  +0.57%  "

Set-TextValue "D39" "0.6305"
Set-TextValue "E39" "  -0.39%  "

Set-TextValue "D40" "4.959"
Set-TextValue "E40" "  +0.66%  "

Set-TextValue "D41" "1.185"
Set-TextValue "E41" "  -1.44%  "

Set-TextValue "D42" "1.390"
Set-TextValue "E42" "  +0.50%  "

Set-TextValue "D43" "7.812"
Set-TextValue "E43" "  +0.85%  "

Set-TextValue "D44" "13.32"
Set-TextValue "E44" "  +0.63%  "

Set-TextValue "D45" "3.726"
Set-TextValue "E45" "  +0.41%  "

Set-TextValue "D46" "0.5893"
Set-TextValue "E46" "  +0.43%  "

Set-TextValue "D47" "122.54"
Set-TextValue "E47" "  +1.17%  "

Set-TextValue "D48" "1.951"
Set-TextValue "E48" "  +0.44%  "

Set-TextValue "D49" "0.06942"
Set-TextValue "E49" "  +1.16%  "

Set-TextValue "D50" "1.135"
Set-TextValue "E50" "  -1.17%  "

Set-TextValue "D51" "72.61"
Set-TextValue "E51" "  +0.72%  "
